$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 557
$ws1.Range("F7").Value = 2773
$ws1.Range("F9").Value = 7779
$ws1.Range("F10").Value = 202
$ws1.Range("F13").Value = 336

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 557
$ws4.Range("F9").Value = 2773
$ws4.Range("F11").Value = 7779
$ws4.Range("F12").Value = 202
$ws4.Range("F17").Value = 336
